# Update the "popularity/views" figures in column F for the "展览" and
# "全部类型" worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 552
$ws1.Range("F4").Value = 197
$ws1.Range("F6").Value = 507
$ws1.Range("F7").Value = 105
$ws1.Range("F10").Value = 6746
$ws1.Range("F11").Value = 233
$ws1.Range("F12").Value = 374
$ws1.Range("F13").Value = 3046
$ws1.Range("F14").Value = 197
$ws1.Range("F15").Value = 347
$ws1.Range("F17").Value = 546

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 552
$ws4.Range("F6").Value = 197
$ws4.Range("F8").Value = 507
$ws4.Range("F9").Value = 105
$ws4.Range("F13").Value = 6746
$ws4.Range("F15").Value = 233
$ws4.Range("F16").Value = 374
$ws4.Range("F17").Value = 3046
$ws4.Range("F18").Value = 197
$ws4.Range("F19").Value = 347
$ws4.Range("F21").Value = 546
